$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A13").Value = "info"
$ws.Range("B13").Value = "file upload"
$ws.Range("C13").Value = "swipl"
$ws.Range("D13").Value = "http://www.swi-prolog.org/howto/http/FileUpload.html"
$ws.Range("D17").Select() | Out-Null
